# Phase4 DVBS2X MODCOD workbook update
#
# The Symbol Rate (B6) drops from 200 to 100 (the units column/labels are
# unaffected). Every downstream throughput figure in column G
# (=E<row>/$B$6) is a live formula, so halving B6 automatically doubles all
# of those cached results on recalculation - no per-cell edits needed there.
#
# The sheet's active-cell selection also moves from C135 back to B5, and the
# book window was repositioned on screen (best effort - harmless if the
# host does not persist raw window geometry).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Symbol Rate (B6): 200 -> 100
$ws.Range("B6").Value = 100

# Recalculate so every dependent formula (the whole G column) gets a fresh
# cached value before save.
$excel.Calculate()

# Restore the active selection to B5
$ws.Range("B5").Select()

# Reposition the workbook window (best effort)
$win = $excel.ActiveWindow
$win.Left = 18100
$win.Top = 1060
